$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.490.10'
$ws.Range("E2").Value = '  -2.65%  '
$ws.Range("D3").Value = '2.470.88'
$ws.Range("E3").Value = '  -2.24%  '
$ws.Range("E4").Value = '  +0.67%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.20'
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '91.37'
$ws.Range("E6").Value = '  -6.93%  '
$ws.Range("E7").Value = '  -3.67%  '
$ws.Range("E8").Value = '  +0.65%  '
$ws.Range("E9").Value = '  -5.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.75'
$ws.Range("E10").Value = '  -6.91%  '
$ws.Range("E11").Value = '  -3.12%  '
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").Value = '2.852.10'
$ws.Range("E13").Value = '  -2.16%  '
$ws.Range("E14").Value = '  -5.47%  '
$ws.Range("D15").Value = '2.541.30'
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.18'
$ws.Range("E16").Value = '  -0.15%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.778'
$ws.Range("E17").Value = '  -4.09%  '
$ws.Range("D18").Value = '41.325.96'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.25'
$ws.Range("E19").Value = '  -5.03%  '
$ws.Range("D20").Value = '0.0₃0915'
$ws.Range("E20").Value = '  -2.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '70.53'
$ws.Range("E21").Value = '  +1.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.92'
$ws.Range("E22").Value = '  -9.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.46'
$ws.Range("E23").Value = '  -3.06%  '
$ws.Range("E24").Value = '  -5.04%  '
$ws.Range("E25").Value = '  +0.15%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.86'
$ws.Range("E26").Value = '  -6.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '23.92'
$ws.Range("E27").Value = '  -6.28%  '
$ws.Range("E28").Value = '  -0.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.65'
$ws.Range("E29").Value = '  -3.45%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.79'
$ws.Range("E30").Value = '  -4.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '152.25'
$ws.Range("E31").Value = '  -2.19%  '
$ws.Range("E32").Value = '  -8.61%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.58'
$ws.Range("E33").Value = '  -5.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.55'
$ws.Range("E34").Value = '  -3.38%  '
$ws.Range("E35").Value = '  -4.03%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.27'
$ws.Range("E36").Value = '  -1.46%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.97'
$ws.Range("E37").Value = '  -5.25%  '
$ws.Range("E38").Value = '  -7.10%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.113'
$ws.Range("E39").Value = '  -3.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0988'
$ws.Range("E40").Value = '  -8.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.02'
$ws.Range("E41").Value = '  -5.25%  '
$ws.Range("E42").Value = '  +0.93%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '1.960.56'
$ws.Range("E43").Value = '  -3.22%  '
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.16'
$ws.Range("E44").Value = '  -10.99%  '
$ws.Range("E45").Value = '  -4.76%  '
$ws.Range("E46").Value = '  -8.72%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.60'
$ws.Range("E47").Value = '  -3.06%  '
$ws.Range("D48").Value = '2.710.23'
$ws.Range("E48").Value = '  -2.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '95.32'
$ws.Range("E49").Value = '  -4.54%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.74'
$ws.Range("E50").Value = '  -5.64%  '
$ws.Range("E51").Value = '  -7.03%  '
